# The deck originally shipped with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" palette (kept only as a leftover /
#                            notes-master-linked part, not the one driving the
#                            visible slide master)
#   ppt/theme/theme2.xml -> "Integral" palette, the theme actually applied to
#                            SlideMaster1 (and therefore every slide's look)
#
# The authored change swaps which palette is "live": the presentation's
# applied design goes from the green/teal "Integral" colour scheme to the
# default blue/orange "Office Theme" colour scheme. We reproduce that by
# rewriting the 12 theme colour slots (dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink) on the presentation's active theme to the stock Office Theme
# values, using the PowerPoint object model's ThemeColorScheme, which is
# exactly the `clrScheme` block inside the theme XML backing the deck's
# slide master.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (scheme slot, target "Office Theme" RGB as BGR-packed integer)
# 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
# 9=accent5 10=accent6 11=hlink 12=folHlink
$tcs.Item(1).RGB  = 0            # dk1      000000
$tcs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388      # dk2      44546A
$tcs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501      # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407        # accent4  FFC000
$tcs.Item(9).RGB  = 12874308     # accent5  4472C4
$tcs.Item(10).RGB = 4697456      # accent6  70AD47
$tcs.Item(11).RGB = 12673797     # hlink    0563C1
$tcs.Item(12).RGB = 7491477      # folHlink 954F72
